$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title heading and matching bold recap line (identical text, both replaced by one global find/replace)
Replace-Text "Play Million 777 for Free - Review of the Online Slot Game" "Play Million 777 for Free - Classic Fruit Machine with Special Features"

# "What we like" bullets
Replace-Text "Offers various special features and bonuses" "Classic fruit machine with 6 reels and 10 paylines"
Replace-Text "Provides an enormous grid with winning possibilities of 1,000,000" "Multiple special features and bonuses to increase chances of winning"
Replace-Text "Traditional and classic graphics that evoke nostalgia among players" "Traditional and classic graphics with a nostalgic feel"
Replace-Text "Offers both low-stake and high-stake betting options" "Accessible betting options for all types of players"

# "What we don't like" bullet
Replace-Text "Not the highest RTP available" "RTP of 95.3% is not the highest available"

# Insert a new bullet paragraph after the RTP bullet, matching its style
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "RTP of 95.3% is not the highest available") {
        $p.Range.InsertParagraphAfter() | Out-Null
        $newPara = $p.Next()
        $newPara.Range.Text = "Limited number of paylines"
        break
    }
}

# Italic summary line
Replace-Text "Read our review of Million 777, a classic fruit machine with exciting special features. Play for free and enjoy the traditional graphics and huge winning possibilities." "Play Million 777 for free and enjoy a classic fruit machine with special features and bonuses."
